$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (column L) held the header "porcentaje-participacion"; it should now
# read "participacion" instead (duplicating the value already used in J2).
$ws.Range("L2").Value = "participacion"

# Row 3 (column L) held "iaest-measure:porcentaje-participacion"; it should
# now read "iaest-measure:participacion" (duplicating J3's value).
$ws.Range("L3").Value = "iaest-measure:participacion"
